{"js": "// Replace the 100 arithmetic-problem cell texts in the single table with\n// their updated values, preserving cell/run formatting. The new values\n// below are the post-edit values in row-major (top-left -> bottom-right)\n// order, matching the table's 20 rows x 5 columns layout.\nconst newValues = [\n  [\"40+56=\", \"9+76=\", \"14+84=\", \"52-42=\", \"85-54=\"],\n  [\"54-23=\", \"0+63=\", \"12+16=\", \"39-34=\", \"95-1=\"],\n  [\"80-8=\", \"20+60=\", \"24-18=\", \"83-21=\", \"11+74=\"],\n  [\"81-25=\", \"5+3=\", \"33+66=\", \"92-32=\", \"20-0=\"],\n  [\"36-6=\", \"80-0=\", \"7-5=\", \"62-39=\", \"84-30=\"],\n  [\"7+32=\", \"7-3=\", \"71-33=\", \"40-8=\", \"59-15=\"],\n  [\"71-12=\", \"7+89=\", \"97-32=\", \"12+1=\", \"73+4=\"],\n  [\"18+3=\", \"50+36=\", \"53+32=\", \"35-33=\", \"60+13=\"],\n  [\"5+71=\", \"91-8=\", \"55+11=\", \"11+55=\", \"12-1=\"],\n  [\"16+45=\", \"35-15=\", \"20+64=\", \"2+51=\", \"35+25=\"],\n  [\"92+4=\", \"90-70=\", \"9+31=\", \"73-57=\", \"43-3=\"],\n  [\"33+32=\", \"45-0=\", \"31+42=\", \"20+72=\", \"83-28=\"],\n  [\"63-35=\", \"51-21=\", \"81-0=\", \"22+59=\", \"49-45=\"],\n  [\"90-51=\", \"46-32=\", \"40-34=\", \"83-26=\", \"42+8=\"],\n  [\"54-30=\", \"83-42=\", \"24+35=\", \"68-67=\", \"86-38=\"],\n  [\"24+58=\", \"10+58=\", \"3+67=\", \"17+1=\", \"92-67=\"],\n  [\"66-39=\", \"1+28=\", \"15+78=\", \"24+59=\", \"17+69=\"],\n  [\"59-6=\", \"80-74=\", \"30-19=\", \"77-9=\", \"72-34=\"],\n  [\"75-32=\", \"4+31=\", \"23+2=\", \"68+2=\", \"6+28=\"],\n  [\"60+35=\", \"15+35=\", \"79-30=\", \"28+68=\", \"66-38=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const oldText = table.values[r][c];\n    const newText = newValues[r][c];\n    if (oldText !== newText) {\n      const cell = table.getCell(r, c);\n      cell.value = newText;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-problem cell texts in the single table with\n# their updated values, preserving cell/run formatting. The new values\n# below are the post-edit values in row-major (top-left -> bottom-right)\n# order, matching the table's 20 rows x 5 columns layout.\n$newValues = @(\n  @(\"40+56=\", \"9+76=\", \"14+84=\", \"52-42=\", \"85-54=\"),\n  @(\"54-23=\", \"0+63=\", \"12+16=\", \"39-34=\", \"95-1=\"),\n  @(\"80-8=\", \"20+60=\", \"24-18=\", \"83-21=\", \"11+74=\"),\n  @(\"81-25=\", \"5+3=\", \"33+66=\", \"92-32=\", \"20-0=\"),\n  @(\"36-6=\", \"80-0=\", \"7-5=\", \"62-39=\", \"84-30=\"),\n  @(\"7+32=\", \"7-3=\", \"71-33=\", \"40-8=\", \"59-15=\"),\n  @(\"71-12=\", \"7+89=\", \"97-32=\", \"12+1=\", \"73+4=\"),\n  @(\"18+3=\", \"50+36=\", \"53+32=\", \"35-33=\", \"60+13=\"),\n  @(\"5+71=\", \"91-8=\", \"55+11=\", \"11+55=\", \"12-1=\"),\n  @(\"16+45=\", \"35-15=\", \"20+64=\", \"2+51=\", \"35+25=\"),\n  @(\"92+4=\", \"90-70=\", \"9+31=\", \"73-57=\", \"43-3=\"),\n  @(\"33+32=\", \"45-0=\", \"31+42=\", \"20+72=\", \"83-28=\"),\n  @(\"63-35=\", \"51-21=\", \"81-0=\", \"22+59=\", \"49-45=\"),\n  @(\"90-51=\", \"46-32=\", \"40-34=\", \"83-26=\", \"42+8=\"),\n  @(\"54-30=\", \"83-42=\", \"24+35=\", \"68-67=\", \"86-38=\"),\n  @(\"24+58=\", \"10+58=\", \"3+67=\", \"17+1=\", \"92-67=\"),\n  @(\"66-39=\", \"1+28=\", \"15+78=\", \"24+59=\", \"17+69=\"),\n  @(\"59-6=\", \"80-74=\", \"30-19=\", \"77-9=\", \"72-34=\"),\n  @(\"75-32=\", \"4+31=\", \"23+2=\", \"68+2=\", \"6+28=\"),\n  @(\"60+35=\", \"15+35=\", \"79-30=\", \"28+68=\", \"66-38=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n  $row = $newValues[$r]\n  for ($c = 0; $c -lt $row.Length; $c++) {\n    $cell = $t.Cell($r + 1, $c + 1)\n    $cell.Range.Text = $row[$c]\n  }\n}\n"}
